# Applies the cryptos-list refresh described by the commit diff
# ("Updated cryptos list ... with GitHub Actions"): refreshed Price /
# Volume(1h) figures throughout, a Litecoin/ShibaInu row swap (rows 17-18),
# and a PaxosStandard -> Quant replacement (row 51) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.388.57'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -3.89%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.846.69'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -5.75%  '
# Row 4
$ws.Range('E4').Value = '  -0.88%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.99'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.17%  '
# Row 6
$ws.Range('E6').Value = '  -0.71%  '
# Row 7
$ws.Range('E7').Value = '  -6.63%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3827'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -5.93%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '49.52'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -7.14%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07815'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -7.67%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.012'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -4.49%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.39'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.56%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.836.21'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -6.69%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.826'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -5.56%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.072'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -7.18%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.005'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.69%  '
# Row 17
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001025'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.49%  '
# Row 18
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '85.02'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -4.88%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06468'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.20%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.94'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -9.04%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.84%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.460'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -6.33%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.368.32'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.03%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.74'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -7.27%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.262'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.20%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.061.26'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -6.59%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '151.27'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.55%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.25'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.80%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.031'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -6.05%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.445'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -8.28%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '119.65'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.30%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.470'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.72%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09278'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.52%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9219'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -6.11%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.591'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.00%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.209'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.90%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02206'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.58%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05921'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -4.99%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.206'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.26%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.271'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -6.33%  '
# Row 41
$ws.Range('E41').Value = '  -0.62%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5884'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.56%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1840'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.22%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.19'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -8.73%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.252'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -6.20%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5637'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -5.57%  '
# Row 47
$ws.Range('E47').Value = '  -6.64%  '
# Row 48
$ws.Range('E48').Value = '  -1.55%  '
# Row 49
$ws.Range('E49').Value = '  -7.14%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06838'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.03%  '
# Row 51
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '107.59'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.73%  '
